# Add I0 and IF columns to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (copy the header style from H1 so I1/J1 match the rest of the header row)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..34: I = 1 (constant), J = same value as H
for ($r = 2; $r -le 34; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
